# Added non active customers who paid in the current month to daily report
# + fixed some bugs.
#
# The "total gviah" (collection) report had a row for every team. A new
# team/customer group ("אודם") now has activity this month, so it is
# inserted at the top of the data (row 2) and every team below shifts its
# displayed name down one row. The last (bottom) row, which used to belong
# to "שנהב", is removed since the data now fits in one fewer row, and all
# the underlying collection figures (columns B/C/D) are refreshed with the
# current month's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New team-name order for rows 2..10 (column A), reflecting the newly
# added team at the top and the corresponding shift of everyone else.
$teams = @(
    "אודם",
    "אלמוג",
    "ברקת",
    "טורקיז",
    "ספיר",
    "פנינה",
    "קריסטל",
    "שוהם - שכר",
    "שנהב"
)

# Updated B (גביה מרגיל) / C (גביה מייעוץ) / D (גביה מפרויקטלי) values for
# each of those rows.
$data = @(
    @(0,     228071, 0),
    @(14420, 315163, 0),
    @(0,     230004, 0),
    @(0,     185964, 0),
    @(0,     0,      0),
    @(0,     1953,   0),
    @(5700,  161796, 0),
    @(0,     17540,  7000),
    @(0,     176550, 0)
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $teams[$i]
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Formula = "=B$row+C$row+D$row"
}

# The report now has one fewer data row (11 -> 10 total rows), so delete
# the old trailing row 11 entirely.
$ws.Rows.Item(11).Delete()
